# Video-Indexing.xlsx — "updated legacy GSC export files"
#
# The GSC (Google Search Console) video-indexing export had two extra
# leading daily rows (2025-09-01 and 2025-09-02) that should be folded
# away: the row for 2025-09-03 becomes the new first data row, but its
# "No video indexed" / "Video indexed" counts are blanked out (the
# Impressions count stays 0) since that day's figures are no longer
# considered valid. Every later row simply shifts up by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the 2025-09-01 row (row 2). Everything below shifts up one row,
# so the 2025-09-02 row is now row 2.
$ws.Rows.Item(2).Delete()

# Remove the (now) 2025-09-02 row. The 2025-09-03 row becomes row 2.
$ws.Rows.Item(2).Delete()

# The surviving 2025-09-03 row keeps its date and Impressions value, but
# "No video indexed" / "Video indexed" are cleared out.
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
